# Updated cryptos list on Mon Jul  1 07:53:28 UTC 2024 with GitHub Actions
# Refreshes the per-coin Price (column D) and Volume(1h) (column E) figures,
# and reflects the row 48/49 swap (dogwifhat <-> InjectiveProtocol) that
# happened in the source ranking between runs.
#
# Column D values are entered with a leading apostrophe so Excel stores them
# as literal text (matching the original inlineStr cells) instead of
# re-parsing "1.00", "7.70", "0.0271", etc. as numbers and normalizing their
# display (which would silently lose trailing zeros / reformat to
# scientific notation). ClearFormats() immediately afterwards drops the
# "number stored as text" quote-prefix styling so the cell is left exactly
# as plain/unstyled as it was before, matching the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.297.90"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "'3.490.72"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "'148.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.75%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.481"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "'7.70"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "'0.126"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").Value = "'4.087.53"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").Value = "'29.84"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.30%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'3.493.40"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "'63.329.36"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "'6.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").Value = "'14.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.39%  "
$ws.Range("D20").Value = "'9.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.66%  "
$ws.Range("D21").Value = "'391.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'0.566"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").Value = "'75.29"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +5.64%  "
$ws.Range("D26").Value = "'3.629.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  -4.48%  "
$ws.Range("D28").Value = "'7.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("D31").Value = "'1.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.33%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'23.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("E35").Value = "  +6.83%  "
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("D37").Value = "'32.33"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +25.80%  "
$ws.Range("D38").Value = "'171.37"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").Value = "'1.57"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.92%  "
$ws.Range("D40").Value = "'3.526.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "'0.809"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.20%  "
$ws.Range("E43").Value = "  +4.50%  "
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "'42.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  +7.91%  "
$ws.Range("D47").Value = "'2.624.84"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.22%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.30"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +13.38%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'23.66"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").Value = "'0.0271"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.54%  "
